$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "skos:prefLabel"
$ws.Range("B9").Value = "EXTRACT"
$ws.Range("C9").Value = "prefLabel of vocabulary"
